$p = $ppt.ActivePresentation

# Add a new slide at the end (position 6) using the "Title and Content" layout
# (ppLayoutText / slideLayout2.xml - matches the layout used by the other
# content slides in this deck).
$s = $p.Slides.Add(6, 2)

# Title placeholder text
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Let’s make it"
